$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/24/2023  Through  7/30/2023"

# --- Row 16 (Burglary): C16 changes type from text "0" to number 1 ---
$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = "#,##0"

# --- Row 17 (Fel. Assault): D17 and E17 change type from number to text ---
#     Copy a same-styled source cell so the shared-string + style (right/center, General) match.
$ws.Range("C15").Copy($ws.Range("D17"))   # -> text "0" (style matches existing s=14 "0" cells)
$ws.Range("H15").Copy($ws.Range("E17"))   # -> text "***.*" (style matches existing s=14 "***.*" cells)

# --- Remaining plain numeric updates ---
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 38
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = -25.490196078431
$ws.Range("L16").Value = 171.428571428571
$ws.Range("M16").Value = -25.490196078431
$ws.Range("N16").Value = -80.104712041884
$ws.Range("C17").Value = 2
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -36.363636363636
$ws.Range("I17").Value = 56
$ws.Range("K17").Value = 12
$ws.Range("L17").Value = 47.368421052631
$ws.Range("M17").Value = 80.645161290322
$ws.Range("N17").Value = -17.647058823529
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 26.666666666666
$ws.Range("I18").Value = 181
$ws.Range("J18").Value = 150
$ws.Range("K18").Value = 20.666666666666
$ws.Range("L18").Value = 34.074074074074
$ws.Range("M18").Value = 36.090225563909
$ws.Range("N18").Value = -69.681742043551
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -27.272727272727
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -22.033898305084
$ws.Range("I19").Value = 381
$ws.Range("J19").Value = 344
$ws.Range("K19").Value = 10.755813953488
$ws.Range("L19").Value = 89.55223880597
$ws.Range("M19").Value = 91.45728643216
$ws.Range("N19").Value = 25.328947368421
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -16.666666666666
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 70
$ws.Range("I20").Value = 91
$ws.Range("J20").Value = 58
$ws.Range("K20").Value = 56.896551724137
$ws.Range("L20").Value = 133.333333333333
$ws.Range("M20").Value = 10.975609756097
$ws.Range("N20").Value = -95.099623047926
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -8.695652173913
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -7.766990291262
$ws.Range("I21").Value = 756
$ws.Range("J21").Value = 655
$ws.Range("K21").Value = 15.419847328244
$ws.Range("L21").Value = 76.223776223776
$ws.Range("M21").Value = 51.2
$ws.Range("N21").Value = -75.016523463317
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -15.384615384615
$ws.Range("F24").Value = 51
$ws.Range("G24").Value = 51
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 332
$ws.Range("J24").Value = 450
$ws.Range("K24").Value = -26.222222222222
$ws.Range("L24").Value = 8.496732026143
$ws.Range("M24").Value = 23.880597014925
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 12
$ws.Range("H25").Value = -20
$ws.Range("I25").Value = 123
$ws.Range("J25").Value = 141
$ws.Range("K25").Value = -12.765957446808
$ws.Range("L25").Value = 55.696202531645
$ws.Range("M25").Value = 23
